$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 - Title shape ("Title 6"): split "And create memories like these.… "
# into four runs: "And create memories like " / "these" / "…" (Marathi-tagged)
# / "… ", mirroring a manual retype + autocorrect edit.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(1)
$tf3 = $shp3.TextFrame

# Neutral placeholder first so PowerPoint's own text-diffing doesn't try to
# keep fragments of the old run around at odd boundaries.
$tf3.TextRange.Text = "XXXXXXXXXXXXXXXXXXXXXXXXXXXXXXXXX"

$tr3 = $tf3.TextRange
$tr3.Text = "… "
$tr3.LanguageID = "en-US"
[void]$tr3.InsertBefore("…")
$tr3.LanguageID = "mr-IN"
[void]$tr3.InsertBefore("these")
$tr3.LanguageID = "en-US"
[void]$tr3.InsertBefore("And create memories like ")
$tr3.LanguageID = "en-US"

# ---------------------------------------------------------------------------
# Slide 4 - "Content Placeholder 5": split the rhetorical question into five
# runs ("fast, " / "easily " / "and " / "conveniently " / "to obtain...") as
# a result of inserting "ly"/"ly" into "easy"/"convenient".
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(4)
$tf4 = $shp4.TextFrame
$tr4 = $tf4.TextRange

$para4 = $tr4.Paragraphs(5, 1)
$para4.Text = "QQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQ"
$para4b = $tr4.Paragraphs(5, 1)
$para4b.Text = "Is the information provided fast, easily and conveniently to obtain and make a quick decision?"

$para4c = $tr4.Paragraphs(5, 1)
$r1 = $para4c.Characters(1, 34)
$r1.Font.Name = "Doctor Soos Bold"
$r2 = $para4c.Characters(35, 7)
$r2.Font.Name = "Doctor Soos Bold"
$r3 = $para4c.Characters(42, 4)
$r3.Font.Name = "Doctor Soos Bold"
$r4 = $para4c.Characters(46, 13)
$r4.Font.Name = "Doctor Soos Bold"
$r5 = $para4c.Characters(59, 36)
$r5.Font.Name = "Doctor Soos Bold"

# ---------------------------------------------------------------------------
# Slide 5 - Title shape ("Title 6"): drop the underline from "The Functionality "
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(1)
$shp5.TextFrame.TextRange.Font.Underline = 0

# ---------------------------------------------------------------------------
# Slide 6 - "Content Placeholder 2": reword the bullet "Increase the data
# provided" -> "Increase location specific information"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(3)
$tf6 = $shp6.TextFrame
$tr6 = $tf6.TextRange

$para6 = $tr6.Paragraphs(5, 1)
$para6.Text = "YYYYYYYYYYYYYYYYYYYYYYYYYYYYYY"
$para6b = $tr6.Paragraphs(5, 1)
$para6b.Text = "Increase location specific information"
